$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the numeric values in row 5 (B5:AH5) to 2 decimal places,
# matching Excel's ROUND() half-away-from-zero behaviour.
$lastCol = 34  # column AH
for ($col = 2; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(5, $col)
    $val = [double]$cell.Value2
    $cell.Value = [Math]::Round($val, 2)
}

# Remove row 6 entirely (data row no longer present after the edit).
$ws.Rows.Item(6).Delete()
